# 2nd May Data Refresh: correct mismatched RegCenter IDs (column A) to 10003
# for the affected device rows, and leave the sheet scrolled/selected at the
# bottom (rows 162 through the end), matching the author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_device")

$rowsToFix = @(3, 23, 43, 63, 83, 105, 114, 123, 132, 141)
foreach ($r in $rowsToFix) {
    $ws.Range("A$r").Value = 10003
}

# Match the saved selection state: rows 162 to the bottom of the sheet selected.
$ws.Rows("162:1048576").Select()
